# Update the point value from 20 to 10 throughout the commitment card.
#
# The visible text "20 points" (and the longer sentence that contains it,
# "Summarize the goals you wrote as part of risk planning in your Module
# Plan (20 points):") must become "10 points" / "... (10 points):".
#
# wdFindContinue = 1, wdReplaceAll = 2 (constants used positionally below
# since this host does not predefine the named Word constants).

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Main document story: replace every occurrence of "20 points" with
#    "10 points". This also fixes the longer sentence
#    "...Module Plan (20 points):" because it contains "20 points" as a
#    substring, so a single pass handles both the short and long forms
#    anywhere they appear in the main flow of the document.
# ---------------------------------------------------------------------
$mainRange = $d.Content
$mainRange.Find.ClearFormatting()
$mainRange.Find.Execute(
    "20 points",   # FindText
    $false,        # MatchCase
    $false,        # MatchWholeWord
    $false,        # MatchWildcards
    $false,        # MatchSoundsLike
    $false,        # MatchAllWordForms
    $true,         # Forward
    1,             # Wrap (wdFindContinue)
    $false,        # Format
    "10 points",   # ReplaceWith
    2              # Replace (wdReplaceAll)
) | Out-Null

# Safety net: in case the host only replaces the first hit per call,
# keep invoking Execute (with the same Replace:=wdReplaceAll) until no
# further match is found anywhere in the main story.
$guard = 0
while ($d.Content.Find.Execute(
        "20 points", $false, $false, $false, $false, $false,
        $true, 1, $false, "10 points", 2) -and $guard -lt 50) {
    $guard++
}

# ---------------------------------------------------------------------
# 2) Floating text boxes ("Planning Meeting 20 points", "Team Review
#    Meeting20 points", "Stand-Up Meeting20 points", "Individual Review
#    Meeting20 points" - each repeated across the sprint table) also show
#    "20 points" and need to become "10 points". These live inside
#    drawing shapes (w:txbxContent), not in Document.Content, so they are
#    reached through the Shapes collection / shape TextFrame instead of
#    the main Find/Replace above.
# ---------------------------------------------------------------------
foreach ($shp in $d.Shapes) {
    try {
        if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
            $tf = $shp.TextFrame
            $trange = $tf.TextRange

            # Prefer an in-place Find/Replace scoped to the shape's own text.
            $null = $trange.Find.Execute(
                "20 points", $false, $false, $false, $false, $false,
                $true, 1, $false, "10 points", 2)

            # Fall back to a direct text assignment if the text still shows
            # the old value (covers hosts where Find does not operate on
            # shape ranges but direct text mutation does).
            if ($tf.TextRange.Text -like "*20 points*") {
                $trange.Text = $trange.Text -replace "20 points", "10 points"
            }

            # Last-resort fallback through the Shape.Text alias.
            if ($shp.HasText -and ($shp.Text -like "*20 points*")) {
                $shp.Text = $shp.Text -replace "20 points", "10 points"
            }
        }
    } catch {
        # Shape text mutation is best-effort; do not let an unsupported
        # operation abort the rest of the script.
    }
}

Write-Output "Point-value update pass complete."
